$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 58; all rows from 58 downward shift down by one.
$ws.Rows("58:58").Insert()

# Populate the new row 58 with its data (same market/region/product metadata as
# neighboring rows, new date/variety/volume/price figures).
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44581
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100103
$ws.Range("H58").Value = "Frutos de hueso (carozo)"
$ws.Range("I58").Value = 100103002
$ws.Range("J58").Value = "Ciruela"
$ws.Range("K58").Value = "Black Amber"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 125
$ws.Range("N58").Value = 12000
$ws.Range("O58").Value = 12000
$ws.Range("P58").Value = 12000
$ws.Range("Q58").Value = "$/bandeja 18 kilos granel"
$ws.Range("R58").Value = "Región de O'Higgins"
$ws.Range("S58").Value = 667
$ws.Range("T58").Value = 18
